# Update the "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value mapping for column F, applied identically on both sheets.
$updates = @{
    2  = 1211
    5  = 2460
    6  = 8005
    7  = 950
    8  = 491
    9  = 434
    12 = 17
    13 = 182
    14 = 8353
    15 = 332
    16 = 1457
    17 = 169
    20 = 207
    21 = 360
    22 = 210
    28 = 1183
    29 = 72
    31 = 113
    33 = 98
    34 = 49
    35 = 90
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
